$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of stock/news data appended as row 36 (date 2025-08-13 / serial 45882)
$ws.Range("A36").Value = 45882
$ws.Range("A36").NumberFormat = $ws.Range("A35").NumberFormat

$ws.Range("B36").Value = 6466.580078125
$ws.Range("C36").Value = 6480.27978515625
$ws.Range("D36").Value = 6445.02001953125
$ws.Range("E36").Value = 6462.669921875
$ws.Range("F36").Value = 3088354000
$ws.Range("G36").Value = 0.0032300788824048
$ws.Range("H36").Value = 1

$i36Text = @"
mp material nysemp stock price up on analyst upgrade affinity bancshares nasdaqafbi v dime community bancshares nasdaqdcom head to head analysis oms energy technology nasdaqomse north american construction group nysenoa headtohead survey reviewing blade air mobility blde and the competition analyzing haverty furniture company nysehvta and brand house collective nasdaqtbhc head to head survey armstrong world industry nyseawi and obayashi otcmktsobycf headtohead review ayro nasdaqayro tesla nasdaqtsla reviewing humana nysehum joint nasdaqjynt financial review kayne anderson bdc nysekbdc hercules capital nysehtgc head to head comparison marui group otcmktsmaury kohl nysekss portillos nasdaqptlo stock rating lowered by wall street zen contrasting bark nysebark and fireman contractor otcmktsfrcn head to head survey southern first bancshares nasdaqsfst v people bancorp of north carolina nasdaqpebk financial comparison city office reit nysecio v lamar advertising nasdaqlamr sadot group nasdaqsdot versus grove collaborative nysegrov headtohead contrast headtohead analysis office property income trust nasdaqopi versus safehold nysesafe financial review ho innovation otcmktsheoff and ho america nasdaqhto lharris technology nyselhx v safran otcmktssafry head to head review orion group nyseorn v mastec nysemtz head to head review headtohead comparison vertiv nysevrt v leidos nyseldos head to head survey hess midstream partner nysehesm v xplr infrastructure nysexifr financial contrast arcutis biotherapeutics nasdaqarqt v novocure nasdaqnvcr contrasting onewater marine nasdaqonew and peloton interactive nasdaqpton critical survey t biosystems nasdaqttoo and perspective therapeutic nysecatx comparing castor maritime nasdaqctrm and united maritime nasdaqusea crypto price today bitcoin price fall to xrp at and solana at riskified ltd nyserskd receives consensus recommendation of hold from brokerage roku inc nasdaqroku given average rating of moderate buy by brokerage financial analysis ameramex international otcmktsammx willis lease finance nasdaqwlfc academy sport and outdoors inc nasdaqaso receives consensus recommendation of moderate buy from analyst brp inc tsedoo given consensus recommendation of hold by brokerage unicredit otcmktsuncry trading down here what happened annaly capital management inc nysenly given average rating of moderate buy by analyst pagseguro digital ltd nysepags receives consensus rating of moderate buy from analyst palantir technology inc nasdaqpltr receives average rating of hold from brokerage travere therapeutic inc nasdaqtvtx given consensus rating of moderate buy by analyst hbo max viu set to offer streaming bundle in asian country sartorius stedim biotech epadim share down should you sell hr real estate investment trust otcmktshruff share up time to buy midnight sun mining cvemma stock price up whats next bmo laddered preferred share index etf tsezpr share down whats next dream industrial real estate investment trust otcmktsdreuf trading down whats next tgs asa otcmktstgsgy trading down whats next doman building material group otcmktscwxzf stock price down should you sell wimi hologram cloud nasdaqwimi trading up still a buy u solar fund lonusfp stock price down whats next titan mining otcmktstnmcf share down whats next china minsheng otcmktscmaky stock price down whats next rational aktiengesellschaft otcmktsrtllf stock price down should you sell intellinetics otcmktsinlxd trading up here what happened melrose industry otcmktsmlspf trading down whats next nippon steel sumitomo metal otcmktsnpscy trading down here what happened barclays etn fi enhanced europe etn series c nysearcaffeu trading up whats next imab nasdaqimab trading higher here why wisdomtree cloud computing fund nasdaqwcld trading down time to sell thai oil public otcmktstoipf trading higher whats next xtrackers international real estate etf nysearcahauz trading down here why fibra uno otcmktsfbasf trading down whats next aimia otcmktsgapff stock price down here what happened camino mineral cvecor share up here what happened sompo otcmktssmpny trading down time to sell foxby otcmktsfxby share up time to buy bmo low volatility canadian equity etf tsezlb trading down time to sell urbana otcmktsubaaf stock price up whats next china citic bank otcmktschcjy stock price down here why biglari nysebha stock price up whats next permrock royalty trust nyseprt trading higher should you buy nufarm otcmktsnufmf trading higher still a buy cascade otcmktscadnf stock price down whats next metaurus u equity exdividend fundseries nysearcaxdiv trading down whats next wisdomtree balanced income fund nysearcawbal share up should you buy midnight sun mining cvemma stock price up whats next uxin nasdaquxin trading up time to buy serica energy otcmktssqzzf trading up whats next bmo laddered preferred share index etf tsezpr share down whats next artis real estate investment trust otcmktsaresf trading down whats next aptus defined risk etf batsdrsk trading down time to sell salmar asa otcmktssalrf share up whats next rpar risk parity etf nysearcarpar share down should you sell sgl carbon otcmktssglff stock price down should you sell overlay share large cap equity etf nysearcaovl share down should you sell boc hong kong otcmktsbhkly share up here what happened amplify alternative harvest etf nysearcamj share up whats next sartorius stedim biotech epadim trading down should you sell manolete partner lonmano stock price down should you sell portfolioplus sp small cap etf nysearcappsc trading higher time to buy vanguard growth etf portfolio tsevgro share down time to sell charlotte web otcmktscwbhf share up here why collplant biotechnology nasdaqclgn trading higher here why iphone pro launch date and price in india revealed check all detail musk grok ai briefly suspended on x for claiming israel commits genocide aersale nasdaqasle raised to strongbuy at wall street zen top tech news nvidias ai chip approved for china meta raise b for ai growth wall street zen downgrade cormedix nasdaqcrmd to buy dragonfly energy dfli expected to announce earnings on thursday mediwound mdwd projected to post quarterly earnings on thursday draftkings nasdaqdkng upgraded at wall street zen palvella therapeutic pvla to release earnings on thursday unusual machine umac expected to announce quarterly earnings on thursday westbow press announces the book answer with truth ha earned exclusive lighthouse recognition portillos nasdaqptlo cut to sell at wall street zen procept biorobotics nasdaqprct rating increased to hold at wall street zen fanatic branded nfl nba mlb nhl apparel sale banned in australia souad m hanna aersale nasdaqasle upgraded at wall street zen donegal group nasdaqdgica stock rating lowered by wall street zen cormedix nasdaqcrmd cut to buy at wall street zen acurx pharmaceutical nasdaqacxp upgraded at wall street zen cricut nasdaqcrct upgraded at wall street zen capital one financial nysecof cut to hold at wall street zen capri nysecpri upgraded at wall street zen citizen nysecia upgraded at wall street zen bwx technology nysebwxt lowered to hold rating by wall street zen vimeo nasdaqvmeo upgraded by wall street zen to buy rating frederick should revisit it district map brehm would be effective leader on council trump say plan to test out putin a europe engages ukraine greg abbott threatens to eliminate democratic district a texas prepares redistricting counterstrike against california map overhaul knightscope nasdaqkscp cut to sell at wall street zen klaviyo nysekvyo upgraded at wall street zen figma nysefig upgraded at wall street zen rxo nyserxo upgraded at wall street zen interactive strength nasdaqtrnr upgraded at wall street zen guardian pharmacy service nysegrdn cut to hold at wall street zen audioeye nasdaqaeye cut to hold at wall street zen biomarin pharmaceutical nasdaqbmrn rating increased to strongbuy at wall street zen banner nasdaqbanr cut to sell at wall street zen alerus financial nasdaqalrs stock rating upgraded by wall street zen ascent solar technology nasdaqasti upgraded at wall street zen ardelyx nasdaqardx upgraded at wall street zen berry nasdaqbry cut to sell at wall street zen applied optoelectronics nasdaqaaoi lowered to sell rating by wall street zen vitesse energy nysevts upgraded at wall street zen amer sport nyseas cut to hold at wall street zen brookdale senior living nysebkd rating increased to hold at wall street zen aflac nyseafl upgraded at wall street zen sapiens international spns expected to announce earnings on wednesday northpointe bancshares inc nysenpb lockup period will expire tomorrow raskin on dc crime do you trust the dems who created the mess or trump who doe the impossible trump back burt jones in republican race for georgia governor in kimmel guest host nicole byer make a dangerous example out of trump dc talk hartselle candidate forum to spotlight local issue u designates baloch separatist a a terror group over role in attack in pakistan eu leader appeal to trump to defend europe security interest at his war summit with putin death notice jean phyllis williams eu leader appeal to trump to defend europe security interest at his war summit with putin the forgotten story of the world deadliest plane crash angler break record after getting green light to reel in harmful fish my fishing pole took a hard hit point coal need to be part of america energy future wdt thought for aug your view constitutional amendment would permit limit on political spending trump back burt jones in republican race for georgia governor in trump back burt jones in republican race for georgia governor in trump back burt jones in republican race for georgia governor in counterpoint america need affordable energy coal doesnt fit the bill u and china extend trade truce another day easing tension between world largest economy u and china extend trade truce another day easing tension between world largest economy at th hour ups bow to teamster strike threat in state nancy pelosi say american know republican ripped away health care will lose midterm unacceptable they want to change rule midgame south korean president lee will travel to washington for aug meeting with trump trump washington police takeover echo history of racist narrative about urban crime hamas hostage video silenced israeli medias talk of gaza aid crisis shareholder action reminder faruqi faruqi llp investigates claim on behalf of investor of lockheed martin fifa move ahead with new human right strategy for world cup game but advocate are skeptical fifa move ahead with new human right strategy for world cup game but advocate are skeptical fifa move ahead with new human right strategy for world cup game but advocate are skeptical inflation likely moved higher last month a tariff bite putting the fed in bind trump washington police takeover echo history of racist narrative about urban crime trump washington police takeover echo history of racist narrative about urban crime trump washington police takeover echo history of racist narrative about urban crime trump washington police takeover echo history of racist narrative about urban crime trump washington police takeover echo history of racist narrative about urban crime federal becomes local the nation capital find itself at the center of a donald trump maelstrom federal becomes local the nation capital find itself at the center of a donald trump maelstrom federal becomes local the nation capital find itself at the center of a donald trump maelstrom letter richard bryan critique donald trump david larson what the dmv can learn from chickfila john hood tariff threat hurt carolina worker letter making my point my turn our right to freedom of conscience and the rowan early voting calendar cartoon what harris campaign memoir should be called editorial hamas steal aid to fund war promote propaganda editorial dallas morning news money role in american politics could get even worse when truth becomes the enemy of the state how to win on immigration go away with leonidas g demas letter stop the hate if you want to up the birth rate scott jennings trigger dems by saying nobody in their right mind think washington dc is safe anna mae midge stellmach trial expose internal tension over trump use of national guard in la photo of mourning for colombian sen miguel uribe who died month after shooting at a rally photo of mourning for colombian sen miguel uribe who died month after shooting at a rally trump unprecedented takeover of dc police come a crime is down figure show application for sd development of affordable housing being accepted library board wont meet in august welcome to south dakota newest news desert with new law city mull election move protest a trump take over dc police deploys national guard in capital dem whistleblower adam schiff okayed leaking of classified info to smear trump in russian collusion hoax letter a reminder from an old friend goodship raise million in series b funding round tesla launch mile model in china a ufc fight at the white house white say it happening a part of deal with paramount cold wallet presale pull in m bonk slip sol breakout loom cenovus indigenous group circle meg energy in possible joint bid pepe coin price prediction at could ozak ai deliver a stronger roi new mexico climate investment center debut energyefficiency loan program for small business nonprofit cardinal health to extend reach by acquiring solaris health dow rally point sp hit new record a cpi report give fed green light to cut rate live update cnbc this best selling jump starter is under for a limited time switzerland could experience unique economic problem from trump tariff bosscattm home service and technology recognized for fourth consecutive year on the inc list of america fastestgrowing private company nuvei tap stablecoins to expand money movement operation fermi americatm name two highly successful nuclear developer mesut and sezin uzman to spearhead nuclear power delivery dow surge point u core inflation rise to meta stock hit new high a thread rocket past million user ai startup perplexity make billion bid for google chrome browser a pause in bitcoin can still be profitable with this option on ibit dutch bros coffee to open th lake houstonarea location in kingwood slqt investor have opportunity to lead selectquote inc security fraud lawsuit with the schall law firm a new chapter for crypto moonshot magax magax bitcoin hyper and aurealone microsoft forecast strong cash flow for fy raising price target to your guide to the top meme coin to invest in trump plan b ipo to privatize fannie mae freddie mac the ally are mostly alright but is that a good thing what to know about trump potential change in federal marijuana policy what to know about trump potential change in federal marijuana policy what to know about trump potential change in federal marijuana policy what to know about trump potential change in federal marijuana policy what to know about trump potential change in federal marijuana policy what to know about trump potential change in federal marijuana policy what to know about trump potential change in federal marijuana policy what to know about trump potential change in federal marijuana policy what to know about trump potential change in federal marijuana policy elon musk threatens apple with lawsuit over openai sparking online feud with sam altman the guardian ibs in america despite advance ibs remains a burden for many million forddirect partner with podium to deliver aipowered conversation to ford and lincoln retailer iconic u film company kodak warns it may go out of business usa today rddt deadline rosen a leading investor right law firm encourages reddit inc investor to crypto entrepreneur who caused market crash pleads guilty to fraud rddt deadline rosen a leading investor right law firm encourages reddit inc investor to secure counsel before important august deadline in security class action rddt gear up for your next project with up to off select tool at the home depot cardinal fiscal q earnings snapshot northern virginia housing market cool slightly a inventory grows resideo technology pop a honeywell indemnity buyout boost growth story iceland eclipse festival to blend music art and space talk hedge fund and insider trading news seth klarman ray dalio ken griffin steve cohen renaissance technology millennium management third point shift payment inc four fastenal co fast and more european auto exporter say tariff stalling business seeing oceanfront property in ukraine trump real estate history shape his presidency spotlight on sezzle analyzing the surge in option activity decoding sofi technologiess option activity whats the big picture peoplesoft how to find own america greatest opportunity why ethereum is rising today rent medical cost drive yearoveryear inflation the magnificent iran and russia vow to block u caucasus ambition best cryptos to buy now in august top crypto opportunity to watch in blockdag xrp stellar solana best dell gaming laptop to buy in abacus insurance and financial service join world insurance associate u july tariff revenue jump fails to halt wider budget deficit bloombergcom best ergonomic gadget for remote worker our view make lasting cut to govguam cost during budget talk why bigbearai stock is cratering today why pubmatic stock fell this morning want to attend world cup game for free fifa need volunteer ai startup perplexity offer to buy google chrome browser for billion the new york time easyfeast unveils new brand identity honoring it mission to make homemade meal easier dow jones leader cisco magnificent seven stock alphabet in new buy zone oryx dental recognized a a top software company on the inc list here why share in united airline took off today advent technology holding inc enhances license to groundbreaking ion pair technology developed at los alamo national laboratory by acquiring exclusivity in marine aviation and portable power field spacex rival ast spacemobile prepares to deploy nearly five dozen satellite signed lanie gardner is on the right track onerpm ink daya g love special sauce generation z saving habit a look at their financial priority and insight zelensky refuse to cede donbas say doing so would give putin springboard for future offensive dorrance planetarium what the cuttingedge led dome at the arizona science center will look like atlantic american corporation report second quarter result for na fallon commander speaks at ceda breakfast atlantic american corporation report second quarter result for gamechange solar rise to third globally and second in u solar tracker market terraform lab founder do kwon pleads guilty to fraud in b terraluna collapse fhlbank san francisco invests million to preserve affordable housing in arizona metro denver inflation hold steady ahead of upcoming tariff keep your toothbrush germfree with the bril uv light toothbrush sanitizer elon musk went from promising tesla unsupervised selfdriving to less nag while los angeles burned rule to protect home from wildfire were on hold spirit airline sound the alarm on it future ability to stay in business what is ai automation and how doe it work official opening of a housing project for adult living with autism spectrum disorder in the gaspe region citizen name barclays executive a it next finance chief worth of labubus headed for black market recovered from socal home los angeles time kucoin first to support ubs umint collateral via digift bridging tokenized rwas and crypto trading bird fan are getting shortchanged reader commentary business spat between daryl hall and john oates ha been resolved in arbitration attorney say desjardins post surplus earnings of million for the second quarter of and surpasses billion in asset family safety stock life surge a parent seek peace of mind bank of america corporation nysebac stock holding decreased by finer wealth management inc tesla inc nasdaqtsla share acquired by saiph capital llc this spy pen secretly record k hd video and audio and it off now the latest national guard arrives in washington dc for trump federal takeover of local police silicon valley lawmaker want to decriminalize homelessness no one is above the law tx ag paxton seek to hold beto orourke in contempt for defying court order ai is reshaping hacking no one agrees how fast trump threatens allowing major lawsuit against fed chair powell over central bank renovation gifford fire grows to more than acre a it shift west in slo county black bear qb embrace life in maine ready to lead football team u teen influencer pilot accused of unauthorized antarctic landing reach deal to leave chile new research discovers the first known latterday saint polygamist with black ancestry israel is in talk to possibly resettle palestinian from gaza in south sudan israel is in talk to possibly resettle palestinian from gaza in south sudan trump takeover of washington law enforcement begin a national guard troop arrive texas senate pass redistricting map despite democrat walkout but house still cant consider bill with no quorum transgender woman sue over jail assault texas republican plan another special session to deliver trump more gop congressional seat texas republican plan another special session to deliver trump more gop congressional seat trump name heritage foundation ej antoni to lead bureau of labor statistic israel is in talk to possibly resettle palestinian from gaza in south sudan israel is in talk to possibly resettle palestinian from gaza in south sudan zelenskyy say putin want the rest of ukraine donetsk region a part of a ceasefire doge land big legal win a appeal court strike down preliminary injunction appeal court remove block on doge data access a familiar name raise democrat spirit in ohio senate race this is how it done here the arrest list from night with trump national guard protecting dc what to know about trump potential change in federal marijuana policy what to know about trump potential change in federal marijuana policy what to know about trump potential change in federal marijuana policy what to know about trump potential change in federal marijuana policy dc top cop clueless on chain of command chief pamela smith deidriven disaster white house say homeless people in dc could be susceptible to jail time white house confirms trumpputin summit in anchorage calling it listening exercise what to know about trump potential change in federal marijuana policy trump use of troop for policing hasnt been seen since america wa ruled by a king trump considering lawsuit against fed a he bash chair jerome powell seeing oceanfront property in ukraine trump real estate history shape his presidency live now u attorney pirro ice and state official announce indictment of haitian gang leader we dont despise the medium enough foreign press association president condemns death of hamas journo how should american college measure merit bank fee for data access threaten the foundation of open banking deja vu another british f just made an emergency landing in japan lima one bill mcdonald honored by housingwire ap news summary at pm edt u citizen detained by ice say she wa denied water for hour israel is in talk to possibly resettle palestinian from gaza in south sudan u designates baloch separatist a a terror group over role in attack in pakistan syria u and jordan say they will work toward a lasting truce in wake of syrian sectarian clash zelensky refuse to cede donbas say doing so would give putin springboard for future offensive paxton urge texas judge to jail beto orourke over fundraising related to redistricting fight more republican join georgia competitive race for lieutenant governor more republican join georgia competitive race for lieutenant governor letter notsofree speech armstrong williams when the bench becomes the battlefield staff commentary takeaway from aps reporting on armed group recruiting child in colombia ashley biden daughter of joe and jill biden file for divorce after year of marriage county call million bond election to build jail live now department of state hold press briefing aug violent guerrilla are taking colombia child unarmed indigenous group are confronting them violent guerrilla are taking colombia child unarmed indigenous group are confronting them boom winsome earlesears just dropped another ad that abigail spanberger will have to answer for watch what to know about the putintrump summit in alaska cdc shooting linked to suspect covid vaccine belief former staff blame rfk jr hostility toward agency the palestinian journalist killed in the war in gaza endured hunger and grief the palestinian journalist killed in the war in gaza endured hunger and grief newsnation cut interview after guest call out trump hypocrisy on live tv newsnation cut interview after guest call out trump hypocrisy on live tv business spat between daryl hall and john oates ha been resolved in arbitration attorney say zelenskyy say putin want the rest of ukraine donetsk region a part of a ceasefire whats at stake a trump prepares to meet putin in alaska shacarri richardson address domestic violence arrest and apologizes to christian coleman musk blow stack after apple snub his racist chatbot texas twostep two additional bill from the texas legislative session that contractor need to know trump want to make the census count a salute to legend like kathy bates jean smart and harrison ford on finally getting their emmy flower inside senate dems summer warpath against trump trump weighs major lawsuit against powell over federal building renovation cost jimmy kimmel obtained italian citizenship in case he need to flee second trump administration it so much worse teen in northfield triple murder to forgo trial and plead guilty teenager in northfield triple murder to plead guilty trump threatens powell with major lawsuit over fed renovation ap news summary at pm edt russian troop advance in ukraine ahead of trumpputin peace summit tinfoil or not you decide damning thread tie nick fuentes to january th in a very troubling way the ncaa is a billiondollar sport behemoth that should not be a nonprofit bosnia serb leader dodik avoids prison term a election loom zelenskyy say russia want a ukrainian withdrawal from ukraine eastern donetsk region a part of a ceasefire deal zelenskyy say russia want a ukrainian withdrawal from ukraine eastern donetsk region a part of a ceasefire deal liar another gop rep face a raucous town hall angered by trump agenda liar another gop rep face a raucous town hall angered by trump agenda poll open in primary election for slain rep melissa hortmans seat pete hegseth belongs to an archconservative church network here what to know cnn slap trump with new nickname accidentally make him sound cooler sherrod brown to launch a comeback attempt for senate in ohio trump v newsom showdown land in court with family tie to the nation highest bench desantis name state sen jay collins a lieutenant governor desantis pick the chuck norris of florida politics a new top deputy allegiant airline to add direct flight from appleton to orlando here what to know bonus editorial cartoon for aug zelenskiy say ukraine wont withdraw troop from donbas first court hearing to take place tuesday in riley strain lawsuit for motion hearing trump washington dc takeover begin a national guard troop arrive
"@
$ws.Range("I36").Value = $i36Text
